$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A115").Value = "TEST"
